$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'232.24"
$ws.Range("D3").Value = "'22.61"
$ws.Range("D4").Value = "'5.300"
$ws.Range("D5").Value = "'0.05608"
$ws.Range("D6").Value = "'3.380"
$ws.Range("D7").Value = "'6.472"
$ws.Range("D8").Value = "'1.061"
$ws.Range("D9").Value = "'0.7822"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1400"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07400"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03177"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02973"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09254"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001670"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.256"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04733"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005791"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006246"
$ws.Range("D20").Value = "'0.005226"
$ws.Range("D21").Value = "'0.001049"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D26").Value = "'0.1300"
$ws.Range("D27").Value = "'0.0004991"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"
$ws.Range("D40").Value = "'0.04033"
$ws.Range("D41").Value = "'0.007012"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1041"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002830"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.01027"
$ws.Range("D45").Value = "'0.00005441"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.7854"
$ws.Range("D48").Value = "'0.04064"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.01010"
